$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.904.94"
$ws.Range("E2").Value = '  -1.69%  '

$ws.Range("D3").Value = "'1.832.32"
$ws.Range("E3").Value = '  -2.00%  '

$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'244.57"
$ws.Range("E5").Value = '  +0.39%  '

$ws.Range("D6").Value = "'0.6899"
$ws.Range("E6").Value = '  -1.96%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = "'0.07687"
$ws.Range("E8").Value = '  -2.99%  '

$ws.Range("E9").Value = '  -2.70%  '

$ws.Range("D10").Value = "'23.42"
$ws.Range("E10").Value = '  -4.37%  '

$ws.Range("D11").Value = "'0.07812"
$ws.Range("E11").Value = '  -0.38%  '

$ws.Range("D12").Value = "'1.835.35"
$ws.Range("E12").Value = '  -3.41%  '

$ws.Range("E13").Value = '  -2.03%  '

$ws.Range("D14").Value = "'90.50"
$ws.Range("E14").Value = '  -3.55%  '

$ws.Range("D15").Value = "'0.6802"
$ws.Range("E15").Value = '  -3.07%  '

$ws.Range("D16").Value = "'6.440"
$ws.Range("E16").Value = '  -1.34%  '

$ws.Range("D17").Value = "'0.000008293"
$ws.Range("E17").Value = '  -1.37%  '

$ws.Range("D18").Value = "'28.941.08"
$ws.Range("E18").Value = '  -2.13%  '

$ws.Range("D19").Value = "'242.90"
$ws.Range("E19").Value = '  -3.87%  '

$ws.Range("D20").Value = "'2.077.87"
$ws.Range("E20").Value = '  -3.13%  '

$ws.Range("E21").Value = '  -2.93%  '

$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("D23").Value = "'7.478"
$ws.Range("E23").Value = '  -2.61%  '

$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").Value = "'162.97"
$ws.Range("E25").Value = '  +0.89%  '

$ws.Range("D26").Value = "'0.1467"
$ws.Range("E26").Value = '  -5.76%  '

$ws.Range("D27").Value = "'8.797"
$ws.Range("E27").Value = '  -2.48%  '

$ws.Range("D28").Value = "'18.20"
$ws.Range("E28").Value = '  -3.42%  '

$ws.Range("D29").Value = "'1.544"
$ws.Range("E29").Value = '  +2.78%  '

$ws.Range("E30").Value = '  -2.62%  '

$ws.Range("D31").Value = "'4.146"
$ws.Range("E31").Value = '  -2.61%  '

$ws.Range("D32").Value = "'1.184"
$ws.Range("E32").Value = '  -2.82%  '

$ws.Range("D33").Value = "'0.05105"
$ws.Range("E33").Value = '  -3.07%  '

$ws.Range("D34").Value = "'0.7656"
$ws.Range("E34").Value = '  +2.16%  '

$ws.Range("D35").Value = "'1.836"
$ws.Range("E35").Value = '  -3.29%  '

$ws.Range("E36").Value = '  -3.36%  '

$ws.Range("D37").Value = "'2.687"
$ws.Range("E37").Value = '  -0.83%  '

$ws.Range("D38").Value = "'0.01847"
$ws.Range("E38").Value = '  -1.80%  '

$ws.Range("D39").Value = "'1.218.95"
$ws.Range("E39").Value = '  -4.35%  '

$ws.Range("E40").Value = '  -2.52%  '

$ws.Range("D41").Value = "'0.9402"
$ws.Range("E41").Value = '  +5.33%  '

$ws.Range("D42").Value = "'108.06"
$ws.Range("E42").Value = '  -1.91%  '

$ws.Range("D43").Value = "'0.9996"
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("D44").Value = "'5.682"
$ws.Range("E44").Value = '  -5.97%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = "'0.00000000122"
$ws.Range("E45").Value = '  -4.10%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'9.552"
$ws.Range("E46").Value = '  -0.82%  '

$ws.Range("D47").Value = "'0.5174"
$ws.Range("E47").Value = '  -0.17%  '

$ws.Range("D48").Value = "'1.978.37"
$ws.Range("E48").Value = '  -3.01%  '

$ws.Range("D49").Value = "'64.17"
$ws.Range("E49").Value = '  -9.73%  '

$ws.Range("D50").Value = "'1.746"
$ws.Range("E50").Value = '  -3.10%  '

$ws.Range("D51").Value = "'0.4187"
$ws.Range("E51").Value = '  -2.70%  '

